$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 4.877553333333334
$ws.Cells.Item(2, 8).Value = 14.63266
$ws.Cells.Item(2, 9).Value = 0.1208715196837975
$ws.Cells.Item(2, 10).Value = 0.1236817629204971
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 5.684044333333333
$ws.Cells.Item(2, 14).Value = 17.052133
$ws.Cells.Item(2, 15).Value = 0.1657316086171094
$ws.Cells.Item(2, 16).Value = 0.1811873991500813
$ws.Cells.Item(2, 17).Value = 27.72422938486444
$ws.Cells.Item(2, 18).Value = 249.51806446378
$ws.Cells.Item(2, 19).Value = 0.02003223139319036
$ws.Cells.Item(2, 20).Value = 0.02240957694586184

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 4.877553333333334
$ws.Cells.Item(3, 8).Value = 14.63266
$ws.Cells.Item(3, 9).Value = 0.1208715196837975
$ws.Cells.Item(3, 10).Value = 0.1236817629204971
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 12.32416133333333
$ws.Cells.Item(3, 14).Value = 36.972484
$ws.Cells.Item(3, 15).Value = 0.3593397522697215
$ws.Cells.Item(3, 16).Value = 0.3928510419240804
$ws.Cells.Item(3, 17).Value = 60.11175419193779
$ws.Cells.Item(3, 18).Value = 541.0057877274401
$ws.Cells.Item(3, 19).Value = 0.04343394193964055
$ws.Cells.Item(3, 20).Value = 0.0485885094303244

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 4.877553333333334
$ws.Cells.Item(4, 8).Value = 14.63266
$ws.Cells.Item(4, 9).Value = 0.1208715196837975
$ws.Cells.Item(4, 10).Value = 0.1236817629204971
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 4.398911333333333
$ws.Cells.Item(4, 14).Value = 13.196734
$ws.Cells.Item(4, 15).Value = 0.1282605498275261
$ws.Cells.Item(4, 16).Value = 0.1402218661287388
$ws.Cells.Item(4, 17).Value = 21.45592463693778
$ws.Cells.Item(4, 18).Value = 193.10332173244
$ws.Cells.Item(4, 19).Value = 0.01550304757313251
$ws.Cells.Item(4, 20).Value = 0.01734288760280436

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 4.877553333333334
$ws.Cells.Item(5, 8).Value = 14.63266
$ws.Cells.Item(5, 9).Value = 0.1208715196837975
$ws.Cells.Item(5, 10).Value = 0.1236817629204971
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 3.112757333333333
$ws.Cells.Item(5, 14).Value = 9.338272
$ws.Cells.Item(5, 15).Value = 0.09075972139462624
$ws.Cells.Item(5, 16).Value = 0.09922378720808876
$ws.Cells.Item(5, 17).Value = 15.18263990705778
$ws.Cells.Item(5, 18).Value = 136.64375916352
$ws.Cells.Item(5, 19).Value = 0.01097026545104654
$ws.Cells.Item(5, 20).Value = 0.01227217292554469

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 4.877553333333334
$ws.Cells.Item(6, 8).Value = 14.63266
$ws.Cells.Item(6, 9).Value = 0.1208715196837975
$ws.Cells.Item(6, 10).Value = 0.1236817629204971
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 13).Value = 8.776807999999999
$ws.Cells.Item(6, 14).Value = 17.553616
$ws.Cells.Item(6, 15).Value = 0.2559083678910167
$ws.Cells.Item(6, 16).Value = 0.1865159055890107
$ws.Cells.Item(6, 17).Value = 42.80934911642667
$ws.Cells.Item(6, 18).Value = 256.85609469856
$ws.Cells.Item(6, 19).Value = 0.03093203332678751
$ws.Cells.Item(6, 20).Value = 0.02306861601596185

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 3.649884
$ws.Cells.Item(7, 8).Value = 10.949652
$ws.Cells.Item(7, 9).Value = 0.09044842682388111
$ws.Cells.Item(7, 10).Value = 0.09255133808384444
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 5.684044333333333
$ws.Cells.Item(7, 14).Value = 17.052133
$ws.Cells.Item(7, 15).Value = 0.1657316086171094
$ws.Cells.Item(7, 16).Value = 0.1811873991500813
$ws.Cells.Item(7, 17).Value = 20.746102467524
$ws.Cells.Item(7, 18).Value = 186.714922207716
$ws.Cells.Item(7, 19).Value = 0.01499016327440873
$ws.Cells.Item(7, 20).Value = 0.01676913623527164

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 3.649884
$ws.Cells.Item(8, 8).Value = 10.949652
$ws.Cells.Item(8, 9).Value = 0.09044842682388111
$ws.Cells.Item(8, 10).Value = 0.09255133808384444
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 12.32416133333333
$ws.Cells.Item(8, 14).Value = 36.972484
$ws.Cells.Item(8, 15).Value = 0.3593397522697215
$ws.Cells.Item(8, 16).Value = 0.3928510419240804
$ws.Cells.Item(8, 17).Value = 44.98175926395201
$ws.Cells.Item(8, 18).Value = 404.835833375568
$ws.Cells.Item(8, 19).Value = 0.03250171528807947
$ws.Cells.Item(8, 20).Value = 0.03635888959770611

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 3.649884
$ws.Cells.Item(9, 8).Value = 10.949652
$ws.Cells.Item(9, 9).Value = 0.09044842682388111
$ws.Cells.Item(9, 10).Value = 0.09255133808384444
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 4.398911333333333
$ws.Cells.Item(9, 14).Value = 13.196734
$ws.Cells.Item(9, 15).Value = 0.1282605498275261
$ws.Cells.Item(9, 16).Value = 0.1402218661287388
$ws.Cells.Item(9, 17).Value = 16.055516092952
$ws.Cells.Item(9, 18).Value = 144.499644836568
$ws.Cells.Item(9, 19).Value = 0.01160096495546575
$ws.Cells.Item(9, 20).Value = 0.01297772133882848

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 3.649884
$ws.Cells.Item(10, 8).Value = 10.949652
$ws.Cells.Item(10, 9).Value = 0.09044842682388111
$ws.Cells.Item(10, 10).Value = 0.09255133808384444
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 3.112757333333333
$ws.Cells.Item(10, 14).Value = 9.338272
$ws.Cells.Item(10, 15).Value = 0.09075972139462624
$ws.Cells.Item(10, 16).Value = 0.09922378720808876
$ws.Cells.Item(10, 17).Value = 11.361203186816
$ws.Cells.Item(10, 18).Value = 102.250828681344
$ws.Cells.Item(10, 19).Value = 0.008209074019117689
$ws.Cells.Item(10, 20).Value = 0.009183294275855262

# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 3.649884
$ws.Cells.Item(11, 8).Value = 10.949652
$ws.Cells.Item(11, 9).Value = 0.09044842682388111
$ws.Cells.Item(11, 10).Value = 0.09255133808384444
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 13).Value = 8.776807999999999
$ws.Cells.Item(11, 14).Value = 17.553616
$ws.Cells.Item(11, 15).Value = 0.2559083678910167
$ws.Cells.Item(11, 16).Value = 0.1865159055890107
$ws.Cells.Item(11, 17).Value = 32.034331090272
$ws.Cells.Item(11, 18).Value = 192.205986541632
$ws.Cells.Item(11, 19).Value = 0.02314650928680947
$ws.Cells.Item(11, 20).Value = 0.01726229663618294

# Row 12
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 16.89730566666667
$ws.Cells.Item(12, 8).Value = 50.691917
$ws.Cells.Item(12, 9).Value = 0.4187351475039348
$ws.Cells.Item(12, 10).Value = 0.428470671797166
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 5.684044333333333
$ws.Cells.Item(12, 14).Value = 17.052133
$ws.Cells.Item(12, 15).Value = 0.1657316086171094
$ws.Cells.Item(12, 16).Value = 0.1811873991500813
$ws.Cells.Item(12, 17).Value = 96.04503452321788
$ws.Cells.Item(12, 18).Value = 864.405310708961
$ws.Cells.Item(12, 19).Value = 0.06939764958034972
$ws.Cells.Item(12, 20).Value = 0.0776334866350166

# Row 13
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 16.89730566666667
$ws.Cells.Item(13, 8).Value = 50.691917
$ws.Cells.Item(13, 9).Value = 0.4187351475039348
$ws.Cells.Item(13, 10).Value = 0.428470671797166
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 12.32416133333333
$ws.Cells.Item(13, 14).Value = 36.972484
$ws.Cells.Item(13, 15).Value = 0.3593397522697215
$ws.Cells.Item(13, 16).Value = 0.3928510419240804
$ws.Cells.Item(13, 17).Value = 208.2451211346476
$ws.Cells.Item(13, 18).Value = 1874.206090211828
$ws.Cells.Item(13, 19).Value = 0.1504681841706892
$ws.Cells.Item(13, 20).Value = 0.1683251498494274

# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 16.89730566666667
$ws.Cells.Item(14, 8).Value = 50.691917
$ws.Cells.Item(14, 9).Value = 0.4187351475039348
$ws.Cells.Item(14, 10).Value = 0.428470671797166
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 4.398911333333333
$ws.Cells.Item(14, 14).Value = 13.196734
$ws.Cells.Item(14, 15).Value = 0.1282605498275261
$ws.Cells.Item(14, 16).Value = 0.1402218661287388
$ws.Cells.Item(14, 17).Value = 74.32974939989757
$ws.Cells.Item(14, 18).Value = 668.9677445990781
$ws.Cells.Item(14, 19).Value = 0.05370720025096491
$ws.Cells.Item(14, 20).Value = 0.060080957180833

# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 16.89730566666667
$ws.Cells.Item(15, 8).Value = 50.691917
$ws.Cells.Item(15, 9).Value = 0.4187351475039348
$ws.Cells.Item(15, 10).Value = 0.428470671797166
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 3.112757333333333
$ws.Cells.Item(15, 14).Value = 9.338272
$ws.Cells.Item(15, 15).Value = 0.09075972139462624
$ws.Cells.Item(15, 16).Value = 0.09922378720808876
$ws.Cells.Item(15, 17).Value = 52.59721212749156
$ws.Cells.Item(15, 18).Value = 473.374909147424
$ws.Cells.Item(15, 19).Value = 0.03800428532559485
$ws.Cells.Item(15, 20).Value = 0.04251448276330883

# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 16.89730566666667
$ws.Cells.Item(16, 8).Value = 50.691917
$ws.Cells.Item(16, 9).Value = 0.4187351475039348
$ws.Cells.Item(16, 10).Value = 0.428470671797166
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 13).Value = 8.776807999999999
$ws.Cells.Item(16, 14).Value = 17.553616
$ws.Cells.Item(16, 15).Value = 0.2559083678910167
$ws.Cells.Item(16, 16).Value = 0.1865159055890107
$ws.Cells.Item(16, 17).Value = 148.3044075536453
$ws.Cells.Item(16, 18).Value = 889.8264453218719
$ws.Cells.Item(16, 19).Value = 0.1071578281763361
$ws.Cells.Item(16, 20).Value = 0.07991659536858019

# Row 17
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 12.177799
$ws.Cells.Item(17, 8).Value = 36.533397
$ws.Cells.Item(17, 9).Value = 0.3017802104744788
$ws.Cells.Item(17, 10).Value = 0.3087965514427589
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 5.684044333333333
$ws.Cells.Item(17, 14).Value = 17.052133
$ws.Cells.Item(17, 15).Value = 0.1657316086171094
$ws.Cells.Item(17, 16).Value = 0.1811873991500813
$ws.Cells.Item(17, 17).Value = 69.21914939842233
$ws.Cells.Item(17, 18).Value = 622.9723445858009
$ws.Cells.Item(17, 19).Value = 0.05001451973074523
$ws.Cells.Item(17, 20).Value = 0.05595004402242777

# Row 18
$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 7).Value = 12.177799
$ws.Cells.Item(18, 8).Value = 36.533397
$ws.Cells.Item(18, 9).Value = 0.3017802104744788
$ws.Cells.Item(18, 10).Value = 0.3087965514427589
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 13).Value = 12.32416133333333
$ws.Cells.Item(18, 14).Value = 36.972484
$ws.Cells.Item(18, 15).Value = 0.3593397522697215
$ws.Cells.Item(18, 16).Value = 0.3928510419240804
$ws.Cells.Item(18, 17).Value = 150.0811595609053
$ws.Cells.Item(18, 18).Value = 1350.730436048148
$ws.Cells.Item(18, 19).Value = 0.1084416260718036
$ws.Cells.Item(18, 20).Value = 0.1213110469768507

# Row 19
$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 7).Value = 12.177799
$ws.Cells.Item(19, 8).Value = 36.533397
$ws.Cells.Item(19, 9).Value = 0.3017802104744788
$ws.Cells.Item(19, 10).Value = 0.3087965514427589
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 13).Value = 4.398911333333333
$ws.Cells.Item(19, 14).Value = 13.196734
$ws.Cells.Item(19, 15).Value = 0.1282605498275261
$ws.Cells.Item(19, 16).Value = 0.1402218661287388
$ws.Cells.Item(19, 17).Value = 53.56905803615533
$ws.Cells.Item(19, 18).Value = 482.1215223253981
$ws.Cells.Item(19, 19).Value = 0.03870649572252319
$ws.Cells.Item(19, 20).Value = 0.04330002869742276

# Row 20
$ws.Cells.Item(20, 5).Value = 3
$ws.Cells.Item(20, 7).Value = 12.177799
$ws.Cells.Item(20, 8).Value = 36.533397
$ws.Cells.Item(20, 9).Value = 0.3017802104744788
$ws.Cells.Item(20, 10).Value = 0.3087965514427589
$ws.Cells.Item(20, 11).Value = 3
$ws.Cells.Item(20, 13).Value = 3.112757333333333
$ws.Cells.Item(20, 14).Value = 9.338272
$ws.Cells.Item(20, 15).Value = 0.09075972139462624
$ws.Cells.Item(20, 16).Value = 0.09922378720808876
$ws.Cells.Item(20, 17).Value = 37.90653314110934
$ws.Cells.Item(20, 18).Value = 341.158798269984
$ws.Cells.Item(20, 19).Value = 0.02738948782507536
$ws.Cells.Item(20, 20).Value = 0.03063996331094795

# Row 21
$ws.Cells.Item(21, 5).Value = 3
$ws.Cells.Item(21, 7).Value = 12.177799
$ws.Cells.Item(21, 8).Value = 36.533397
$ws.Cells.Item(21, 9).Value = 0.3017802104744788
$ws.Cells.Item(21, 10).Value = 0.3087965514427589
$ws.Cells.Item(21, 11).Value = 2
$ws.Cells.Item(21, 13).Value = 8.776807999999999
$ws.Cells.Item(21, 14).Value = 17.553616
$ws.Cells.Item(21, 15).Value = 0.2559083678910167
$ws.Cells.Item(21, 16).Value = 0.1865159055890107
$ws.Cells.Item(21, 17).Value = 106.882203685592
$ws.Cells.Item(21, 18).Value = 641.2932221135519
$ws.Cells.Item(21, 19).Value = 0.07722808112433135
$ws.Cells.Item(21, 20).Value = 0.0575954684351097

# Row 22
$ws.Cells.Item(22, 5).Value = 2
$ws.Cells.Item(22, 7).Value = 2.750664
$ws.Cells.Item(22, 8).Value = 5.501328
$ws.Cells.Item(22, 9).Value = 0.06816469551390787
$ws.Cells.Item(22, 10).Value = 0.04649967575573358
$ws.Cells.Item(22, 11).Value = 3
$ws.Cells.Item(22, 13).Value = 5.684044333333333
$ws.Cells.Item(22, 14).Value = 17.052133
$ws.Cells.Item(22, 15).Value = 0.1657316086171094
$ws.Cells.Item(22, 16).Value = 0.1811873991500813
$ws.Cells.Item(22, 17).Value = 15.634896122104
$ws.Cells.Item(22, 18).Value = 93.80937673262399
$ws.Cells.Item(22, 19).Value = 0.01129704463841542
$ws.Cells.Item(22, 20).Value = 0.00842515531150346

# Row 23
$ws.Cells.Item(23, 5).Value = 2
$ws.Cells.Item(23, 7).Value = 2.750664
$ws.Cells.Item(23, 8).Value = 5.501328
$ws.Cells.Item(23, 9).Value = 0.06816469551390787
$ws.Cells.Item(23, 10).Value = 0.04649967575573358
$ws.Cells.Item(23, 11).Value = 3
$ws.Cells.Item(23, 13).Value = 12.32416133333333
$ws.Cells.Item(23, 14).Value = 36.972484
$ws.Cells.Item(23, 15).Value = 0.3593397522697215
$ws.Cells.Item(23, 16).Value = 0.3928510419240804
$ws.Cells.Item(23, 17).Value = 33.899626909792
$ws.Cells.Item(23, 18).Value = 203.397761458752
$ws.Cells.Item(23, 19).Value = 0.02449428479950865
$ws.Cells.Item(23, 20).Value = 0.01826744606977184

# Row 24
$ws.Cells.Item(24, 5).Value = 2
$ws.Cells.Item(24, 7).Value = 2.750664
$ws.Cells.Item(24, 8).Value = 5.501328
$ws.Cells.Item(24, 9).Value = 0.06816469551390787
$ws.Cells.Item(24, 10).Value = 0.04649967575573358
$ws.Cells.Item(24, 11).Value = 3
$ws.Cells.Item(24, 13).Value = 4.398911333333333
$ws.Cells.Item(24, 14).Value = 13.196734
$ws.Cells.Item(24, 15).Value = 0.1282605498275261
$ws.Cells.Item(24, 16).Value = 0.1402218661287388
$ws.Cells.Item(24, 17).Value = 12.099927043792
$ws.Cells.Item(24, 18).Value = 72.599562262752
$ws.Cells.Item(24, 19).Value = 0.008742841325439723
$ws.Cells.Item(24, 20).Value = 0.006520271308850237

# Row 25
$ws.Cells.Item(25, 5).Value = 2
$ws.Cells.Item(25, 7).Value = 2.750664
$ws.Cells.Item(25, 8).Value = 5.501328
$ws.Cells.Item(25, 9).Value = 0.06816469551390787
$ws.Cells.Item(25, 10).Value = 0.04649967575573358
$ws.Cells.Item(25, 11).Value = 3
$ws.Cells.Item(25, 13).Value = 3.112757333333333
$ws.Cells.Item(25, 14).Value = 9.338272
$ws.Cells.Item(25, 15).Value = 0.09075972139462624
$ws.Cells.Item(25, 16).Value = 0.09922378720808876
$ws.Cells.Item(25, 17).Value = 8.562149537536
$ws.Cells.Item(25, 18).Value = 51.372897225216
$ws.Cells.Item(25, 19).Value = 0.006186608773791807
$ws.Cells.Item(25, 20).Value = 0.004613873932432033

# Row 26
$ws.Cells.Item(26, 5).Value = 2
$ws.Cells.Item(26, 7).Value = 2.750664
$ws.Cells.Item(26, 8).Value = 5.501328
$ws.Cells.Item(26, 9).Value = 0.06816469551390787
$ws.Cells.Item(26, 10).Value = 0.04649967575573358
$ws.Cells.Item(26, 11).Value = 2
$ws.Cells.Item(26, 13).Value = 8.776807999999999
$ws.Cells.Item(26, 14).Value = 17.553616
$ws.Cells.Item(26, 15).Value = 0.2559083678910167
$ws.Cells.Item(26, 16).Value = 0.1865159055890107
$ws.Cells.Item(26, 17).Value = 24.142049800512
$ws.Cells.Item(26, 18).Value = 96.56819920204799
$ws.Cells.Item(26, 19).Value = 0.01744391597675227
$ws.Cells.Item(26, 20).Value = 0.008672929133176015
